# "before demo on create entity&sendnotify"
#
# Refresh the sample phone-number values used by the CreateEntity / SendNotify
# demo data on Sheet1-Sheet4 (MobileNumber, Enquiry_PhoneNumber, Lead_PN,
# Sales_PN in row 2). Sheet5 stores these as plain numbers and is untouched.
#
# The numbers must stay text cells (they were shared strings before the
# edit), but a bare numeric-looking string like "9840062226" gets
# auto-coerced to a number by Excel. Prefixing with an apostrophe forces
# Excel to store it as text without requiring a new/changed number format,
# so re-normalising the style afterwards restores the original look
# (General format, top-aligned) without minting a new cell style.

$wb = $excel.ActiveWorkbook

$xlTop = -4160

$newValues = @{
    "G2"  = "9840062226"
    "AF2" = "9840039741"
    "AV2" = "9840072901"
    "AZ2" = "9840014664"
}

$sheetNames = @("Sheet1", "Sheet2", "Sheet3", "Sheet4")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    foreach ($addr in $newValues.Keys) {
        $cell = $ws.Range($addr)
        $cell.Value = "'" + $newValues[$addr]
        $cell.Style = "Normal"
        $cell.VerticalAlignment = $xlTop
    }
}
